# Populate the "ibbs-web-app-metrics" sheet with the new performance-test
# scenario rows (rows 2-21) captured for the "IBBS APEX app".
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ibbs-web-app-metrics")

$ws1.Cells.Item(2, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(2, 2).Value = "Local"
$ws1.Cells.Item(2, 3).Value = "Local"
$ws1.Cells.Item(2, 4).Value = 45560.930810185186
$ws1.Cells.Item(2, 5).Value = 45560.51414351852
$ws1.Cells.Item(2, 6).Value = "Login"
$ws1.Cells.Item(2, 7).Value = "Page Load"
$ws1.Cells.Item(2, 8).Value = 25
$ws1.Cells.Item(2, 9).Value = 1682.94
$ws1.Cells.Item(2, 10).Value = 5.772
$ws1.Cells.Item(2, 11).Value = "Login.png"

$ws1.Cells.Item(3, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(3, 2).Value = "Local"
$ws1.Cells.Item(3, 3).Value = "Local"
$ws1.Cells.Item(3, 4).Value = 45560.9308912037
$ws1.Cells.Item(3, 5).Value = 45560.51422453704
$ws1.Cells.Item(3, 6).Value = "Full Sampling Plan Report"
$ws1.Cells.Item(3, 7).Value = "Login/Page Load"
$ws1.Cells.Item(3, 8).Value = 41
$ws1.Cells.Item(3, 9).Value = 3799.84
$ws1.Cells.Item(3, 10).Value = 3.376
$ws1.Cells.Item(3, 11).Value = "Full Sampling Plan Report.png"

$ws1.Cells.Item(4, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(4, 2).Value = "Local"
$ws1.Cells.Item(4, 3).Value = "Local"
$ws1.Cells.Item(4, 4).Value = 45560.9309375
$ws1.Cells.Item(4, 5).Value = 45560.51427083334
$ws1.Cells.Item(4, 6).Value = "Full Sampling Plan Report"
$ws1.Cells.Item(4, 7).Value = "Page Reload/Filter Report"
$ws1.Cells.Item(4, 8).Value = 41
$ws1.Cells.Item(4, 9).Value = 3800.15
$ws1.Cells.Item(4, 10).Value = 2.241
$ws1.Cells.Item(4, 11).Value = "Full Sampling Plan Report filter.png"

$ws1.Cells.Item(5, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(5, 2).Value = "Local"
$ws1.Cells.Item(5, 3).Value = "Local"
$ws1.Cells.Item(5, 4).Value = 45560.93096064815
$ws1.Cells.Item(5, 5).Value = 45560.51429398148
$ws1.Cells.Item(5, 6).Value = "View Specimens"
$ws1.Cells.Item(5, 7).Value = "Page Load"
$ws1.Cells.Item(5, 8).Value = 34
$ws1.Cells.Item(5, 9).Value = 4538.41
$ws1.Cells.Item(5, 10).Value = 5.8
$ws1.Cells.Item(5, 11).Value = "View Specimens.png"

$ws1.Cells.Item(6, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(6, 2).Value = "Local"
$ws1.Cells.Item(6, 3).Value = "Local"
$ws1.Cells.Item(6, 4).Value = 45560.93104166666
$ws1.Cells.Item(6, 5).Value = 45560.514375
$ws1.Cells.Item(6, 6).Value = "View/Edit Specimen"
$ws1.Cells.Item(6, 7).Value = "Page Load"
$ws1.Cells.Item(6, 8).Value = 35
$ws1.Cells.Item(6, 9).Value = 3131.02
$ws1.Cells.Item(6, 10).Value = 1.787
$ws1.Cells.Item(6, 11).Value = "View Edit Specimen.png"

$ws1.Cells.Item(7, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(7, 2).Value = "Local"
$ws1.Cells.Item(7, 3).Value = "Local"
$ws1.Cells.Item(7, 4).Value = 45560.931180555555
$ws1.Cells.Item(7, 5).Value = 45560.51451388889
$ws1.Cells.Item(7, 6).Value = "View/Edit Specimen"
$ws1.Cells.Item(7, 7).Value = "Form submission"
$ws1.Cells.Item(7, 8).Value = 54
$ws1.Cells.Item(7, 9).Value = 3827.29
$ws1.Cells.Item(7, 10).Value = 2.4
$ws1.Cells.Item(7, 11).Value = "View Edit Specimen post specimen record insert.png"

$ws1.Cells.Item(8, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(8, 2).Value = "Local"
$ws1.Cells.Item(8, 3).Value = "Local"
$ws1.Cells.Item(8, 4).Value = 45560.93126157407
$ws1.Cells.Item(8, 5).Value = 45560.51459490741
$ws1.Cells.Item(8, 6).Value = "View Specimens"
$ws1.Cells.Item(8, 7).Value = "Form submission"
$ws1.Cells.Item(8, 8).Value = 40
$ws1.Cells.Item(8, 9).Value = 4546.86
$ws1.Cells.Item(8, 10).Value = 5.636
$ws1.Cells.Item(8, 11).Value = "View Specimens post specimen record update.png"

$ws1.Cells.Item(9, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(9, 2).Value = "Local"
$ws1.Cells.Item(9, 3).Value = "Local"
$ws1.Cells.Item(9, 4).Value = 45560.93133101852
$ws1.Cells.Item(9, 5).Value = 45560.51466435185
$ws1.Cells.Item(9, 6).Value = "Download Specimen Data"
$ws1.Cells.Item(9, 7).Value = "IBBS_SPEC_DATA_YYYYMMDD.csv"
$ws1.Cells.Item(9, 8).Value = 1
$ws1.Cells.Item(9, 9).Value = 1811934
$ws1.Cells.Item(9, 10).Value = 2.532
$ws1.Cells.Item(9, 11).Value = "View Specimens specimen download complete.png"

$ws1.Cells.Item(10, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(10, 2).Value = "Local"
$ws1.Cells.Item(10, 3).Value = "Local"
$ws1.Cells.Item(10, 4).Value = 45560.93136574074
$ws1.Cells.Item(10, 5).Value = 45560.514699074076
$ws1.Cells.Item(10, 6).Value = "Sampling Plan Summary Region Report"
$ws1.Cells.Item(10, 7).Value = "Page Load"
$ws1.Cells.Item(10, 8).Value = 47
$ws1.Cells.Item(10, 9).Value = 5405.12
$ws1.Cells.Item(10, 10).Value = 5.419
$ws1.Cells.Item(10, 11).Value = "Sampling Plan Summary Region Report.png"

$ws1.Cells.Item(11, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(11, 2).Value = "Local"
$ws1.Cells.Item(11, 3).Value = "Local"
$ws1.Cells.Item(11, 4).Value = 45560.931435185186
$ws1.Cells.Item(11, 5).Value = 45560.51476851852
$ws1.Cells.Item(11, 6).Value = "Sampling Plan Summary Region Report"
$ws1.Cells.Item(11, 7).Value = "Page Reload/Filter Report"
$ws1.Cells.Item(11, 8).Value = 46
$ws1.Cells.Item(11, 9).Value = 5168.35
$ws1.Cells.Item(11, 10).Value = 4.526
$ws1.Cells.Item(11, 11).Value = "Sampling Plan Summary Region Report filter.png"

$ws1.Cells.Item(12, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(12, 2).Value = "Local"
$ws1.Cells.Item(12, 3).Value = "Remote"
$ws1.Cells.Item(12, 4).Value = 45560.93430555556
$ws1.Cells.Item(12, 5).Value = 45560.51763888889
$ws1.Cells.Item(12, 6).Value = "Login"
$ws1.Cells.Item(12, 7).Value = "Page Load"
$ws1.Cells.Item(12, 8).Value = 26
$ws1.Cells.Item(12, 9).Value = 745.47
$ws1.Cells.Item(12, 10).Value = 2.737
$ws1.Cells.Item(12, 11).Value = "Login.png"

$ws1.Cells.Item(13, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(13, 2).Value = "Local"
$ws1.Cells.Item(13, 3).Value = "Remote"
$ws1.Cells.Item(13, 4).Value = 45560.93435185185
$ws1.Cells.Item(13, 5).Value = 45560.51768518519
$ws1.Cells.Item(13, 6).Value = "Full Sampling Plan Report"
$ws1.Cells.Item(13, 7).Value = "Login/Page Load"
$ws1.Cells.Item(13, 8).Value = 40
$ws1.Cells.Item(13, 9).Value = 1158.84
$ws1.Cells.Item(13, 10).Value = 3.675
$ws1.Cells.Item(13, 11).Value = "Full Sampling Plan Report.png"

$ws1.Cells.Item(14, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(14, 2).Value = "Local"
$ws1.Cells.Item(14, 3).Value = "Remote"
$ws1.Cells.Item(14, 4).Value = 45560.93439814815
$ws1.Cells.Item(14, 5).Value = 45560.51773148148
$ws1.Cells.Item(14, 6).Value = "Full Sampling Plan Report"
$ws1.Cells.Item(14, 7).Value = "Page Reload/Filter Report"
$ws1.Cells.Item(14, 8).Value = 40
$ws1.Cells.Item(14, 9).Value = 1158.86
$ws1.Cells.Item(14, 10).Value = 2.225
$ws1.Cells.Item(14, 11).Value = "Full Sampling Plan Report filter.png"

$ws1.Cells.Item(15, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(15, 2).Value = "Local"
$ws1.Cells.Item(15, 3).Value = "Remote"
$ws1.Cells.Item(15, 4).Value = 45560.934432870374
$ws1.Cells.Item(15, 5).Value = 45560.5177662037
$ws1.Cells.Item(15, 6).Value = "View Specimens"
$ws1.Cells.Item(15, 7).Value = "Page Load"
$ws1.Cells.Item(15, 8).Value = 35
$ws1.Cells.Item(15, 9).Value = 1391.29
$ws1.Cells.Item(15, 10).Value = 7.115
$ws1.Cells.Item(15, 11).Value = "View Specimens.png"

$ws1.Cells.Item(16, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(16, 2).Value = "Local"
$ws1.Cells.Item(16, 3).Value = "Remote"
$ws1.Cells.Item(16, 4).Value = 45560.934525462966
$ws1.Cells.Item(16, 5).Value = 45560.517858796295
$ws1.Cells.Item(16, 6).Value = "View/Edit Specimen"
$ws1.Cells.Item(16, 7).Value = "Page Load"
$ws1.Cells.Item(16, 8).Value = 33
$ws1.Cells.Item(16, 9).Value = 911.34
$ws1.Cells.Item(16, 10).Value = 3.799
$ws1.Cells.Item(16, 11).Value = "View Edit Specimen.png"

$ws1.Cells.Item(17, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(17, 2).Value = "Local"
$ws1.Cells.Item(17, 3).Value = "Remote"
$ws1.Cells.Item(17, 4).Value = 45560.934699074074
$ws1.Cells.Item(17, 5).Value = 45560.51803240741
$ws1.Cells.Item(17, 6).Value = "View/Edit Specimen"
$ws1.Cells.Item(17, 7).Value = "Form submission"
$ws1.Cells.Item(17, 8).Value = 53
$ws1.Cells.Item(17, 9).Value = 1174.57
$ws1.Cells.Item(17, 10).Value = 4.884
$ws1.Cells.Item(17, 11).Value = "View Edit Specimen post specimen record insert.png"

$ws1.Cells.Item(18, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(18, 2).Value = "Local"
$ws1.Cells.Item(18, 3).Value = "Remote"
$ws1.Cells.Item(18, 4).Value = 45560.934907407405
$ws1.Cells.Item(18, 5).Value = 45560.51824074074
$ws1.Cells.Item(18, 6).Value = "View Specimens"
$ws1.Cells.Item(18, 7).Value = "Form submission"
$ws1.Cells.Item(18, 8).Value = 41
$ws1.Cells.Item(18, 9).Value = 1392.18
$ws1.Cells.Item(18, 10).Value = 6.832
$ws1.Cells.Item(18, 11).Value = "View Specimens post specimen record update.png"

$ws1.Cells.Item(19, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(19, 2).Value = "Local"
$ws1.Cells.Item(19, 3).Value = "Remote"
$ws1.Cells.Item(19, 4).Value = 45560.935
$ws1.Cells.Item(19, 5).Value = 45560.51833333333
$ws1.Cells.Item(19, 6).Value = "Download Specimen Data"
$ws1.Cells.Item(19, 7).Value = "IBBS_SPEC_DATA_YYYYMMDD.csv"
$ws1.Cells.Item(19, 8).Value = 1
$ws1.Cells.Item(19, 9).Value = 1804315
$ws1.Cells.Item(19, 10).Value = 3.162
$ws1.Cells.Item(19, 11).Value = "View Specimens specimen download complete.png"

$ws1.Cells.Item(20, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(20, 2).Value = "Local"
$ws1.Cells.Item(20, 3).Value = "Remote"
$ws1.Cells.Item(20, 4).Value = 45560.9350462963
$ws1.Cells.Item(20, 5).Value = 45560.518379629626
$ws1.Cells.Item(20, 6).Value = "Sampling Plan Summary Region Report"
$ws1.Cells.Item(20, 7).Value = "Page Load"
$ws1.Cells.Item(20, 8).Value = 44
$ws1.Cells.Item(20, 9).Value = 1512.01
$ws1.Cells.Item(20, 10).Value = 5.207
$ws1.Cells.Item(20, 11).Value = "Sampling Plan Summary Region Report.png"

$ws1.Cells.Item(21, 1).Value = "IBBS APEX app"
$ws1.Cells.Item(21, 2).Value = "Local"
$ws1.Cells.Item(21, 3).Value = "Remote"
$ws1.Cells.Item(21, 4).Value = 45560.93513888889
$ws1.Cells.Item(21, 5).Value = 45560.518472222226
$ws1.Cells.Item(21, 6).Value = "Sampling Plan Summary Region Report"
$ws1.Cells.Item(21, 7).Value = "Page Reload/Filter Report"
$ws1.Cells.Item(21, 8).Value = 48
$ws1.Cells.Item(21, 9).Value = 1516.82
$ws1.Cells.Item(21, 10).Value = 5.506
$ws1.Cells.Item(21, 11).Value = "Sampling Plan Summary Region Report filter.png"

# The worksheet's formulas recalc automatically once the COM script
# finishes, so the "Summary" sheet's aggregate statistics update to
# reflect the newly-entered rows without any further action here.

# View-state updates to match what the workbook looked like when the
# author finished entering this data: the data-entry range got selected
# on the metrics sheet, and focus moved over to the Summary sheet.
$ws1.Range("A1:K21").Select()

$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Activate()
$ws2.Range("F9").Select()
